$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-05 21:48:03"
$ws.Range("E3").Value = "2026-02-05 21:48:06"
$ws.Range("E4").Value = "2026-02-05 21:48:09"
$ws.Range("O4").Value = "11.7 °C"
$ws.Range("E5").Value = "2026-02-05 21:48:11"
$ws.Range("E6").Value = "2026-02-05 21:48:14"
$ws.Range("E7").Value = "2026-02-05 21:48:17"
$ws.Range("H7").Value = "'78%"
$ws.Range("E8").Value = "2026-02-05 21:48:19"
$ws.Range("O8").Value = "9.1 °C"
$ws.Range("E9").Value = "2026-02-05 21:48:22"
$ws.Range("E10").Value = "2026-02-05 21:48:25"
$ws.Range("E11").Value = "2026-02-05 21:48:28"
$ws.Range("J11").Value = "994.5 hPa"
$ws.Range("O11").Value = "0.9 °C"
$ws.Range("E12").Value = "2026-02-05 21:48:30"
$ws.Range("H12").Value = "'84%"
$ws.Range("O12").Value = "10.5 °C"
$ws.Range("E13").Value = "2026-02-05 21:48:33"
$ws.Range("E14").Value = "2026-02-05 21:48:35"
$ws.Range("I14").Value = "8.1 mm"
$ws.Range("E15").Value = "2026-02-05 21:48:38"
$ws.Range("H15").Value = "'79%"
$ws.Range("O15").Value = "8.9 °C"
$ws.Range("E16").Value = "2026-02-05 21:48:41"
$ws.Range("E17").Value = "2026-02-05 21:48:43"
$ws.Range("I17").Value = "8.7 mm"
$ws.Range("M17").Value = "2.7 °C 21:18 TU"
$ws.Range("E18").Value = "2026-02-05 21:48:46"
$ws.Range("E19").Value = "2026-02-05 21:48:49"
$ws.Range("E20").Value = "2026-02-05 21:48:52"
$ws.Range("E21").Value = "2026-02-05 21:48:54"
$ws.Range("H21").Value = "'81%"
$ws.Range("J21").Value = "990.5 hPa"
$ws.Range("O21").Value = "6.6 °C"
$ws.Range("E22").Value = "2026-02-05 21:48:57"
$ws.Range("H22").Value = "'85%"
$ws.Range("O22").Value = "9.3 °C"
$ws.Range("E23").Value = "2026-02-05 21:49:00"
$ws.Range("H23").Value = "'87%"
$ws.Range("O23").Value = "8.1 °C"
$ws.Range("E24").Value = "2026-02-05 21:49:02"
$ws.Range("E25").Value = "2026-02-05 21:49:05"
$ws.Range("J25").Value = "993.9 hPa"
$ws.Range("E26").Value = "2026-02-05 21:49:08"
$ws.Range("E27").Value = "2026-02-05 21:49:11"
$ws.Range("H27").Value = "'90%"
$ws.Range("E28").Value = "2026-02-05 21:49:13"
$ws.Range("J28").Value = "992.7 hPa"
$ws.Range("O28").Value = "2.9 °C"
$ws.Range("E29").Value = "2026-02-05 21:49:16"
$ws.Range("O29").Value = "9.6 °C"
$ws.Range("E30").Value = "2026-02-05 21:49:19"
$ws.Range("E31").Value = "2026-02-05 21:49:21"
$ws.Range("I31").Value = "19.6 mm"
$ws.Range("E32").Value = "2026-02-05 21:49:24"
$ws.Range("H32").Value = "'79%"
$ws.Range("E33").Value = "2026-02-05 21:49:27"
$ws.Range("E34").Value = "2026-02-05 21:49:29"
$ws.Range("H34").Value = "'94%"
$ws.Range("O34").Value = "4.5 °C"
$ws.Range("E35").Value = "2026-02-05 21:49:32"
$ws.Range("E36").Value = "2026-02-05 21:49:35"
